# "fixed formatting of java code"
#
# 1) Slide 14, shape "Rechteck 7" (the second Java code box, the
#    "public class Contact { ... }" snippet): resize/move it slightly
#    (left edge moves left, width grows) and bump the indent level of
#    the 6 "private ..." field-declaration paragraphs from level 0 to
#    level 1 (PowerPoint's 1-based IndentLevel 1 -> 2).
# 2) Slide 17, shape "Textfeld 1": nudge its left edge slightly to the
#    right.
#
# NOTE: this COM runtime stores Shape.Left/Top/Width/Height internally
# as 32-bit (single-precision) floats expressed in points, and floors
# the point->EMU conversion (EMU = floor(float32(points) * 12700)).
# To land exactly on the target EMU values from the OOXML diff we feed
# in points values that have been solved so that this float32+floor
# round trip reproduces the exact target EMU, rather than the naive
# emu/12700 conversion (which can be off by one EMU).

$p = $ppt.ActivePresentation

# --- Slide 14: resize "Rechteck 7" and re-indent the field lines ---
$s14 = $p.Slides.Item(14)
$rechteck7 = $s14.Shapes.Item("Rechteck 7")

$rechteck7.Left = 548.3738708677165   # -> a:off x  = 6964348 EMU
$rechteck7.Width = 388.8422699645669  # -> a:ext cx = 4938297 EMU

$tr = $rechteck7.TextFrame.TextRange
# Paragraphs 2..7 are the "private String <field>;" lines; bump them to
# indent level 1 (0-based OOXML lvl="1"), i.e. COM IndentLevel 2.
for ($i = 2; $i -le 7; $i++) {
    $tr.Paragraphs($i, 1).IndentLevel = 2
}

# --- Slide 17: nudge "Textfeld 1" to the right ---
$s17 = $p.Slides.Item(17)
$textfeld1 = $s17.Shapes.Item("Textfeld 1")
$textfeld1.Left = 86.30047234094489   # -> a:off x = 1096016 EMU
